# Revert the "Finance - Core Banking System Modernization" KPI dashboard
# content back from the generic "Artificial Intelligence and Machine
# Learning" template text, restoring the original Finance-flavored KPI
# rows, values and formulas.

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions & User Guide")
$wsDashboard    = $wb.Worksheets.Item("KPI Dashboard")

# ---------------------------------------------------------------------
# Sheet: Instructions & User Guide
# ---------------------------------------------------------------------
$wsInstructions.Range("A1").Value = "Finance - Core Banking System Modernization KPI Dashboard - User Guide & Instructions"

# ---------------------------------------------------------------------
# Sheet: KPI Dashboard
# ---------------------------------------------------------------------
$wsDashboard.Range("A1").Value = "FINANCE - CORE BANKING SYSTEM MODERNIZATION - KPI DASHBOARD"
$wsDashboard.Range("A2").Value = "Project: Core Banking System Modernization"

# Note text shared by every KPI row
$noteText = "Critical KPI for Finance - Core Banking System Modernization success"

# KPI001
$wsDashboard.Range("B8").Value = "Transaction Processing Accuracy"
$wsDashboard.Range("C8").Value = "Performance"
$wsDashboard.Range("D8").Value = 99.95
$wsDashboard.Range("E8").Value = 99.87
$wsDashboard.Range("F8").Formula = "=((E8-D8)/D8)*100"
$wsDashboard.Range("G8").Value = "At Risk"
$wsDashboard.Range("K8").Value = $noteText

# KPI002
$wsDashboard.Range("B9").Value = "System Uptime %"
$wsDashboard.Range("C9").Value = "Performance"
$wsDashboard.Range("D9").Value = 99.9
$wsDashboard.Range("E9").Value = 99.92
$wsDashboard.Range("F9").Formula = "=((E9-D9)/D9)*100"
$wsDashboard.Range("G9").Value = "On Track"
$wsDashboard.Range("K9").Value = $noteText

# KPI003
$wsDashboard.Range("B10").Value = "Data Migration Completion"
$wsDashboard.Range("C10").Value = "Performance"
$wsDashboard.Range("D10").Value = 95
$wsDashboard.Range("E10").Value = 88
$wsDashboard.Range("F10").Formula = "=((E10-D10)/D10)*100"
$wsDashboard.Range("G10").Value = "At Risk"
$wsDashboard.Range("K10").Value = $noteText

# KPI004
$wsDashboard.Range("B11").Value = "User Adoption Rate"
$wsDashboard.Range("C11").Value = "Performance"
$wsDashboard.Range("D11").Value = 85
$wsDashboard.Range("E11").Value = 78
$wsDashboard.Range("F11").Formula = "=((E11-D11)/D11)*100"
$wsDashboard.Range("G11").Value = "At Risk"
$wsDashboard.Range("K11").Value = $noteText

# KPI005
$wsDashboard.Range("B12").Value = "Regulatory Compliance Score"
$wsDashboard.Range("C12").Value = "Quality"
$wsDashboard.Range("D12").Value = 100
$wsDashboard.Range("E12").Value = 98
$wsDashboard.Range("F12").Formula = "=((E12-D12)/D12)*100"
$wsDashboard.Range("G12").Value = "On Track"
$wsDashboard.Range("K12").Value = $noteText

# KPI006
$wsDashboard.Range("B13").Value = "Security Incident Count"
$wsDashboard.Range("C13").Value = "Quality"
$wsDashboard.Range("D13").Value = 0
$wsDashboard.Range("E13").Value = 2
$wsDashboard.Range("F13").Formula = "=((E13-D13)/D13)*100"
$wsDashboard.Range("G13").Value = "At Risk"
$wsDashboard.Range("K13").Value = $noteText

# KPI007
$wsDashboard.Range("B14").Value = "Defect Density"
$wsDashboard.Range("C14").Value = "Quality"
$wsDashboard.Range("D14").Value = 0.5
$wsDashboard.Range("E14").Value = 0.8
$wsDashboard.Range("F14").Formula = "=((E14-D14)/D14)*100"
$wsDashboard.Range("G14").Value = "At Risk"
$wsDashboard.Range("K14").Value = $noteText

# KPI008
$wsDashboard.Range("B15").Value = "Training Completion Rate"
$wsDashboard.Range("C15").Value = "Quality"
$wsDashboard.Range("D15").Value = 95
$wsDashboard.Range("E15").Value = 92
$wsDashboard.Range("F15").Formula = "=((E15-D15)/D15)*100"
$wsDashboard.Range("G15").Value = "On Track"
$wsDashboard.Range("K15").Value = $noteText

# KPI009
$wsDashboard.Range("B16").Value = "Budget Variance %"
$wsDashboard.Range("C16").Value = "Financial"
$wsDashboard.Range("D16").Value = 0
$wsDashboard.Range("E16").Value = 3.5
$wsDashboard.Range("F16").Formula = "=((E16-D16)/D16)*100"
$wsDashboard.Range("G16").Value = "At Risk"
$wsDashboard.Range("K16").Value = $noteText

# KPI010
$wsDashboard.Range("B17").Value = "Schedule Variance %"
$wsDashboard.Range("C17").Value = "Financial"
$wsDashboard.Range("D17").Value = 0
$wsDashboard.Range("E17").Value = 2.8
$wsDashboard.Range("F17").Formula = "=((E17-D17)/D17)*100"
$wsDashboard.Range("G17").Value = "At Risk"
$wsDashboard.Range("K17").Value = $noteText

# KPI011
$wsDashboard.Range("B18").Value = "Cost Savings Achieved"
$wsDashboard.Range("C18").Value = "Financial"
$wsDashboard.Range("D18").Value = 2500000
$wsDashboard.Range("E18").Value = 2100000
$wsDashboard.Range("F18").Formula = "=((E18-D18)/D18)*100"
$wsDashboard.Range("G18").Value = "At Risk"
$wsDashboard.Range("K18").Value = $noteText

# KPI012
$wsDashboard.Range("B19").Value = "ROI Achievement %"
$wsDashboard.Range("C19").Value = "Financial"
$wsDashboard.Range("D19").Value = 125
$wsDashboard.Range("E19").Value = 118
$wsDashboard.Range("F19").Formula = "=((E19-D19)/D19)*100"
$wsDashboard.Range("G19").Value = "At Risk"
$wsDashboard.Range("K19").Value = $noteText

# KPI013
$wsDashboard.Range("B20").Value = "Customer Satisfaction Score"
$wsDashboard.Range("C20").Value = "Financial"
$wsDashboard.Range("D20").Value = 4.5
$wsDashboard.Range("E20").Value = 4.2
$wsDashboard.Range("F20").Formula = "=((E20-D20)/D20)*100"
$wsDashboard.Range("G20").Value = "On Track"
$wsDashboard.Range("K20").Value = $noteText

# KPI014
$wsDashboard.Range("B21").Value = "Integration Success Rate"
$wsDashboard.Range("C21").Value = "Financial"
$wsDashboard.Range("D21").Value = 98
$wsDashboard.Range("E21").Value = 95
$wsDashboard.Range("F21").Formula = "=((E21-D21)/D21)*100"
$wsDashboard.Range("G21").Value = "On Track"
$wsDashboard.Range("K21").Value = $noteText

# KPI015
$wsDashboard.Range("C22").Value = "Financial"
$wsDashboard.Range("D22").Value = 100
$wsDashboard.Range("E22").Value = 92
$wsDashboard.Range("F22").Formula = "=((E22-D22)/D22)*100"
$wsDashboard.Range("G22").Value = "At Risk"
$wsDashboard.Range("K22").Value = $noteText
